$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the COMPLIANCE_GROUPS column (column I) entirely - the remote
# access sheet now only tracks SERVER_GROUPS.
$ws.Columns.Item(9).Delete()

# Header: GROUPS -> SERVER_GROUPS
$ws.Range("H1").Value = "SERVER_GROUPS"

# Row 2: fill in the host that previously had no HOST value (the cell
# already carries the row's formatting).
$ws.Range("A2").Value = "server01.example.com"

# Row 3 keeps its host/type/node/key/user but its row height grows a lot
# (the SSH key cell wraps across many more lines in the new layout).
$ws.Rows.Item(3).RowHeight = 1572.35

# Row 4: replace the WinRM slave entry with another SSH server entry
# using the new SERVER_GROUPS column (cyberwatch) instead of WinRM creds.
$ws.Range("A4").Style = $ws.Range("F4").Style
$ws.Range("A4").Value = "10.0.2.16"

$ws.Range("B4").Value = 22

$ws.Range("C4").Value = "CbwRam::RemoteAccess::Ssh::WithPassword"

# D4 goes back to the plain/default look the rest of the NODE column uses.
$ws.Range("D4").Clear()
$ws.Range("D4").Value = "master"

# E4 stays empty but keeps the "group" cell formatting like E2.
$ws.Range("E4").Style = $ws.Range("F4").Style

$ws.Range("F4").Value = "cyberwatch"
$ws.Range("G4").Value = "cyberwatch"
$ws.Range("H4").Value = "production, test"

# Header row is a touch taller in the new layout.
$ws.Rows.Item(1).RowHeight = 14.9

# Selection moves to the newly-filled row 2 (A2 active, whole row selected).
$ws.Range("A2").Select() | Out-Null
$ws.Rows.Item(2).Select() | Out-Null
